$d = $word.ActiveDocument

# 1. "As an example" -> "As a relevant example"
$d.Content.Find.Execute("As an example of anchoring", $true, $false, $false, $false, $false, $true, 1, $false, "As a relevant example of anchoring", 2)

# 2. "in 2015 may" -> "in March may"
$d.Content.Find.Execute("in 2015 may induce buyers in 2016 to pay", $true, $false, $false, $false, $false, $true, 1, $false, "in March may induce buyers in April to pay", 2)

Write-Host "done"
